# Fixes for target values
$wb = $excel.ActiveWorkbook

# Work on the "ModeloA" sheet (sheet1.xml) - append 3 new rows (Jul, Ago, Set)
$ws = $wb.Worksheets.Item("ModeloA")

$ws.Cells.Item(8, 1).Value = "Jul"
$ws.Cells.Item(8, 2).Value = 10000
$ws.Cells.Item(8, 3).Value = 11000

$ws.Cells.Item(9, 1).Value = "Ago"
$ws.Cells.Item(9, 2).Value = 11000
$ws.Cells.Item(9, 3).Value = 11000

$ws.Cells.Item(10, 1).Value = "Set"
$ws.Cells.Item(10, 2).Value = 12000
$ws.Cells.Item(10, 3).Value = 11000

# Update selection on ModeloA to C1 and make it the active/selected tab
$ws.Range("C1").Select()
$ws.Activate()

# ReceitasM (sheet8) was previously the selected tab; move selection
# there back off of "tabSelected" by activating ModeloA above and
# leaving ReceitasM's own in-sheet selection untouched.
